$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Prefix "G-" onto each existing product_code value in column G (rows 2-14).
# Row 6 is intentionally blank and must remain blank.
for ($r = 2; $r -le 14; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2
    if ($val -ne $null -and $val -ne "") {
        $cell.Value = "G-" + $val
    }
}

# Update the active selection to G15, matching the post-edit saved state.
$ws.Range("G15").Select()
